# Update SYDATA - Work order test cases on the "RMA Details Maintenance Grid" sheet
# Replaces the RMA-3ZTQ-* sample data (rows 2-4, columns E/F/J) with the new
# RMA-PPLQ-* sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-PPLQ-001"
$ws.Range("F2").Value = "RMA-PPLQ-1-1"
$ws.Range("J2").Value = "a7s5f000000xL3IAAU"

# Row 3
$ws.Range("E3").Value = "RMA-PPLQ-002"
$ws.Range("F3").Value = "RMA-PPLQ-1-2"
$ws.Range("J3").Value = "a7s5f000000xL3JAAU"

# Row 4
$ws.Range("E4").Value = "RMA-PPLQ-003"
$ws.Range("F4").Value = "RMA-PPLQ-1-3"
$ws.Range("J4").Value = "a7s5f000000xL3KAAU"

# Columns F (RMA line number) and J (Salesforce Id) are best-fit / autofit,
# so re-run AutoFit after the content change to refresh their widths.
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(10).AutoFit()
